# RPA Calendar v2 수정
# Replace the (September) calendar grid in Sheet1 with the new (June) calendar
# grid, update the associated row heights, and mark the workbook window as
# minimized, matching the authoritative diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 1 (height unchanged: 43.5) ---
$ws.Range("A1").Value = 28
$ws.Range("B1").Value = 29
$ws.Range("C1").Value = 30
$ws.Range("D1").Value = 31
$ws.Range("E1").Value = "1`n의병의 날"
$ws.Range("F1").Value = 2
$ws.Range("G1").Value = "3`n음4.15"

# --- Row 2 (height 101.5 -> 72.5) ---
$ws.Range("A2").Value = 4
$ws.Range("B2").Value = "5`n환경의 날"
$ws.Range("C2").Value = "6`n현충일`n망종"
$ws.Range("D2").Value = 7
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 9
$ws.Range("G2").Value = "10`n6.10 민주항쟁기념일"

# --- Row 3 (height 43.5 -> 72.5) ---
$ws.Range("A3").Value = 11
$ws.Range("B3").Value = 12
$ws.Range("C3").Value = 13
$ws.Range("D3").Value = 14
$ws.Range("E3").Value = "15`n노인학대 예방의 날"
$ws.Range("F3").Value = 16
$ws.Range("G3").Value = 17

# --- Row 4 (height 72.5 -> 58) ---
$ws.Range("A4").Value = "18`n음5.1"
$ws.Range("B4").Value = 19
$ws.Range("C4").Value = 20
$ws.Range("D4").Value = "21`n하지`n해양조사의 날"
$ws.Range("E4").Value = "22`n단오"
$ws.Range("F4").Value = 23
$ws.Range("G4").Value = 24

# --- Row 5 (height default -> 43.5) ---
$ws.Range("A5").Value = "25`n6·25 전쟁일"
$ws.Range("B5").Value = 26
$ws.Range("C5").Value = 27
$ws.Range("D5").Value = "28`n철도의 날"
$ws.Range("E5").Value = 29
$ws.Range("F5").Value = 30
$ws.Range("G5").Value = 1

# Row 6 is unchanged.

# --- Row heights ---
$ws.Rows.Item(1).RowHeight = 43.5
$ws.Rows.Item(2).RowHeight = 72.5
$ws.Rows.Item(3).RowHeight = 72.5
$ws.Rows.Item(4).RowHeight = 58
$ws.Rows.Item(5).RowHeight = 43.5
$ws.Rows.Item(6).RowHeight = 87

# --- Mark the workbook window as minimized ---
$excel.Windows.Item(1).WindowState = -4140

$wb.Save()
